$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '64.240.66'
$ws.Range('E2').Value2 = '  -4.02%  '
$ws.Range('D3').Value2 = '3.014.99'
$ws.Range('E3').Value2 = '  -6.16%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value2 = '0.998'
$ws.Range('E4').Value2 = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '553.94'
$ws.Range('E5').Value2 = '  -6.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '141.53'
$ws.Range('E6').Value2 = '  -6.33%  '
$ws.Range('E7').Value2 = '  -0.28%  '
$ws.Range('D8').Value2 = '2.995.99'
$ws.Range('E8').Value2 = '  -6.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.475'
$ws.Range('E9').Value2 = '  -12.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '0.152'
$ws.Range('E10').Value2 = '  -12.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '6.11'
$ws.Range('E11').Value2 = '  -6.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '0.446'
$ws.Range('E12').Value2 = '  -10.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '33.98'
$ws.Range('E13').Value2 = '  -12.77%  '
$ws.Range('E14').Value2 = '  -13.91%  '
$ws.Range('D15').Value2 = '3.485.98'
$ws.Range('E15').Value2 = '  -6.64%  '
$ws.Range('D16').Value2 = '64.142.36'
$ws.Range('E16').Value2 = '  -4.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '0.110'
$ws.Range('E17').Value2 = '  -4.16%  '
$ws.Range('D18').Value2 = '3.005.51'
$ws.Range('E18').Value2 = '  -6.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '479.49'
$ws.Range('E19').Value2 = '  -9.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '6.37'
$ws.Range('E20').Value2 = '  -11.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '13.27'
$ws.Range('E21').Value2 = '  -11.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '0.647'
$ws.Range('E22').Value2 = '  -14.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '6.71'
$ws.Range('E23').Value2 = '  -15.57%  '
$ws.Range('B24').Value2 = 'Litecoin'
$ws.Range('C24').Value2 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '77.45'
$ws.Range('E24').Value2 = '  -9.52%  '
$ws.Range('B25').Value2 = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '12.28'
$ws.Range('E25').Value2 = '  -11.12%  '
$ws.Range('E26').Value2 = '  +0.24%  '
$ws.Range('E27').Value2 = '  -15.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value2 = '7.69'
$ws.Range('E28').Value2 = '  -5.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '2.01'
$ws.Range('E29').Value2 = '  -8.93%  '
$ws.Range('B30').Value2 = 'Stacks'
$ws.Range('C30').Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value2 = '2.56'
$ws.Range('E30').Value2 = '  -3.87%  '
$ws.Range('B31').Value2 = 'EthereumClassic'
$ws.Range('C31').Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '25.33'
$ws.Range('E31').Value2 = '  -13.52%  '
$ws.Range('B32').Value2 = 'Mantle'
$ws.Range('C32').Value2 = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '1.12'
$ws.Range('E32').Value2 = '  -1.36%  '
$ws.Range('B33').Value2 = 'Bittensor'
$ws.Range('C33').Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '514.41'
$ws.Range('E33').Value2 = '  -5.63%  '
$ws.Range('B34').Value2 = 'FirstDigitalUSD'
$ws.Range('C34').Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '0.999'
$ws.Range('E34').Value2 = '  -0.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '51.90'
$ws.Range('E35').Value2 = '  -3.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '5.06'
$ws.Range('E36').Value2 = '  -11.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '5.70'
$ws.Range('E37').Value2 = '  -12.60%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '0.0407'
$ws.Range('E38').Value2 = '  -4.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '0.122'
$ws.Range('E39').Value2 = '  -3.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '0.0774'
$ws.Range('E40').Value2 = '  -10.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '8.05'
$ws.Range('E41').Value2 = '  -13.83%  '
$ws.Range('D42').Value2 = '2.765.42'
$ws.Range('E42').Value2 = '  -5.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '2.41'
$ws.Range('E43').Value2 = '  -8.60%  '
$ws.Range('E44').Value2 = '  -0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '0.233'
$ws.Range('E45').Value2 = '  -11.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '1.98'
$ws.Range('E46').Value2 = '  -6.32%  '
$ws.Range('B47').Value2 = 'Monero'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '115.40'
$ws.Range('E47').Value2 = '  -6.89%  '
$ws.Range('B48').Value2 = 'PEPE'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D48').Value2 = '0.0₃0505'
$ws.Range('E48').Value2 = '  -13.37%  '
$ws.Range('E49').Value2 = '  -9.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '23.18'
$ws.Range('E50').Value2 = '  -12.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '1.99'
$ws.Range('E51').Value2 = '  -17.76%  '
